$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.174721837043762
$ws.Range("B1").Value = 2.398706436157227
$ws.Range("D1").Value = 2.355192184448242
$ws.Range("E1").Value = 1.208419919013977
